$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Diretores de departamento" to "juri"
$ws.Name = "juri"
